$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.787.95'
$ws.Range("E2").Value = '  -3.96%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.817.34'
$ws.Range("E3").Value = '  -2.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '277.47'
$ws.Range("E5").Value = '  -7.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5093'
$ws.Range("E7").Value = '  -4.69%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3518'
$ws.Range("E8").Value = '  -6.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.36'
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06654'
$ws.Range("E10").Value = '  -7.21%  '
$ws.Range("E11").Value = '  -7.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.8331'
$ws.Range("E12").Value = '  -6.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07903'
$ws.Range("E13").Value = '  -2.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.826.63'
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.081'
$ws.Range("E15").Value = '  -3.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.63'
$ws.Range("E16").Value = '  -5.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.12'
$ws.Range("E18").Value = '  -4.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008021'
$ws.Range("E19").Value = '  -6.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9999'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '25.842.59'
$ws.Range("E21").Value = '  -3.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.724'
$ws.Range("E22").Value = '  -4.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.01'
$ws.Range("E23").Value = '  -6.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.080'
$ws.Range("E24").Value = '  -4.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.44'
$ws.Range("E25").Value = '  -3.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.182'
$ws.Range("E26").Value = '  -3.19%  '
$ws.Range("E27").Value = '  -3.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.11'
$ws.Range("E28").Value = '  -4.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '109.57'
$ws.Range("E29").Value = '  -4.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.334'
$ws.Range("E30").Value = '  -8.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.241'
$ws.Range("E31").Value = '  -7.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08832'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04865'
$ws.Range("E33").Value = '  -2.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7337'
$ws.Range("E34").Value = '  -8.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.136'
$ws.Range("E35").Value = '  -2.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.878'
$ws.Range("E36").Value = '  -3.87%  '
$ws.Range("E37").Value = '  +0.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9995'
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5234'
$ws.Range("E39").Value = '  -12.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.319'
$ws.Range("E40").Value = '  -11.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01847'
$ws.Range("E41").Value = '  -5.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9609'
$ws.Range("E42").Value = '  -10.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '111.65'
$ws.Range("E43").Value = '  -3.52%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.199'
$ws.Range("E44").Value = '  -6.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.058'
$ws.Range("E45").Value = '  -9.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9999'
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4615'
$ws.Range("E47").Value = '  -9.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1363'
$ws.Range("E48").Value = '  -8.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.50'
$ws.Range("E49").Value = '  -3.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.187'
$ws.Range("E50").Value = '  -7.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.501'
$ws.Range("E51").Value = '  -7.55%  '
